$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@('GUGGU UPI_504H Tws Bluetoth for all Smart phones Blueto...', '₹859')
    ,@('SYARA UPI_504H_ K1 Bluetooth Headset for all Smart phon...', '₹379')
    ,@('ROAR UPI_504H K1 Bluetooth Headset', '₹449')
    ,@('GUGGU UPI_504H Bluetooth Headset', '₹439')
    ,@('ROAR UPI_504H_ K1 Bluetooth Headset for all Smart phone...', '₹499')
    ,@('SYARA UPI_504H Kaju Bluetooth Headset for all Smart pho...', '₹439')
    ,@('GUGGU UPI_504H K1 Bluetoth Headset for all Smart phones...', '₹379')
    ,@('ROAR UPI_504H_ K1 Bluetooth Headset for all Smart phone...', '₹449')
    ,@('SYARA UPI_504H KAJU Blutooth Hadst for all Smart phons ...', '₹439')
    ,@('Jio V2 Keypad Phone (UPI) Locked', '₹1,056')
    ,@('Jio B1 Keypad Phone (UPI) Locked Blue', '₹1,419')
    ,@('Jio B1 Keypad Phone (UPI) Locked', '₹1,400')
    ,@('Jio V2 Keypad Phone with (UPI) Locked New', '₹1,068')
    ,@('mrquee upi hot Heating pad for pain relief, heating bag...', '₹399')
    ,@('SYARA UPI_504H HSB-730 Bluetooth Headset for all Smartp...', '₹459')
    ,@('SYARA UPI_504HI 7 Twins Bluetooth Headset for all Smart...', '₹859')
    ,@('SYARA UPI_504H TWS4 Bluetooth Headset for all Smart pho...', '₹959')
    ,@('ROAR UPI_504H_Bluetooth Headset for all Smart phones Bl...', '₹749')
    ,@('ROAR UPI_504H_mi Bluetooth Headset', '₹449')
    ,@('SYARA UPI_504H K1 Bluetoth Headset for all Smart phones...', '₹379')
    ,@('SYARA UPI_504H_m-i TWS-4 Bluetoth Headset for all Smart...', '₹799')
    ,@('SYARA UPI_504HK1 Bluetooth Headset for all Smart phones...', '₹379')
    ,@('GUGGU UPI_504H Ms771 Bluetooth for all Smart phones Blu...', '₹999')
    ,@('itel SG400 Keypad Mobile Phone with 9.8mm Ultra Slim De...', '₹1,375')
    ,@('itel SG400 Keypad Mobile Phone with 9.8mm Ultra Slim De...', '₹1,398')
    ,@('itel SG200 Keypad Phone with 1200mAh Battery 1.3 MP Cam...', '₹1,190')
    ,@('itel SG200 Keypad Phone with 1200mAh Battery 1.3 MP Cam...', '₹1,190')
    ,@('HEPOCH Screen Guard for Fire-Boltt Phoenix Smart Watch ...', '₹189')
    ,@('GUGGU UPI_504H I7 Twins Bluetooth for all Smart phones ...', '₹449')
    ,@('SYARA UPI_504H_Bluetooth Headset for all Smart phones B...', '₹499')
    ,@('SYARA UPI_504H L21 Bluetooth Headset for all Smart phon...', '₹1,089')
    ,@('Beautum Model EBMSU022482 MISS YOU Upi Name Printed Bes...', '₹344')
    ,@('badri creations Happy Birthday Upi White Coffee Ceramic...', '₹279')
    ,@('GUGGU UPI_324E_M4 Fitness band', '₹559')
    ,@('GUGGU UPI_324E_M3 Fitness band', '₹449')
    ,@('Clairbell UPI_177Y_A1 Smart Watch Smartwatch', '₹1,109')
    ,@('SYARA UPI_264L_ Smart Watch Smartwatch', '₹1,059')
    ,@('itel SG200 Keypad Phone with 1200mAh Battery 1.3 MP Cam...', '₹1,329')
    ,@('itel SG400 Keypad Mobile Phone with 9.8mm Ultra Slim De...', '₹1,376')
    ,@('SYARA UPI_564L YT228 Portable tripod Compatiable with a...', '₹1,419')
    ,@($null, '₹1,190')
    ,@($null, '₹1,398')
    ,@($null, '₹1,400')
    ,@($null, '₹379')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $a = $data[$i][0]
    $b = $data[$i][1]
    if ($a -ne $null) {
        $ws.Cells.Item($r, 1).Value = $a
    } else {
        $ws.Cells.Item($r, 1).Value = ""
    }
    if ($b -ne $null) {
        $ws.Cells.Item($r, 2).Value = $b
    }
}
